$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.122.44'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.163.32'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '602.36'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '153.90'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.548'
$ws.Range("E8").Value = '  +2.68%  '
$ws.Range("D9").Value = '3.161.20'
$ws.Range("E9").Value = '  -1.72%  '
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("E11").Value = '  -10.99%  '
$ws.Range("D12").Value = '0.516'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = '0.0000266'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").Value = '38.32'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '3.678.73'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").Value = '66.190.44'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '7.41'
$ws.Range("E17").Value = '  -0.91%  '
$ws.Range("D18").Value = '3.156.16'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").Value = '510.15'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("D21").Value = '15.41'
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").Value = '8.05'
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '14.59'
$ws.Range("E24").Value = '  -4.50%  '
$ws.Range("D25").Value = '84.40'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = '2.99'
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = '9.06'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").Value = '2.39'
$ws.Range("E29").Value = '  +6.22%  '
$ws.Range("D30").Value = '3.06'
$ws.Range("E30").Value = '  +6.04%  '
$ws.Range("D31").Value = '6.93'
$ws.Range("E31").Value = '  +0.67%  '
$ws.Range("D32").Value = '27.98'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("E35").Value = '  -2.24%  '
$ws.Range("D36").Value = '500.54'
$ws.Range("E36").Value = '  +3.62%  '
$ws.Range("D37").Value = '54.69'
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("E38").Value = '  -3.25%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").Value = '0.129'
$ws.Range("E40").Value = '  +8.22%  '
$ws.Range("D41").Value = '8.74'
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").Value = '0.0₃0674'
$ws.Range("E42").Value = '  +4.81%  '
$ws.Range("D43").Value = '0.295'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  -6.73%  '
$ws.Range("D45").Value = '2.42'
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D46").Value = '2.825.00'
$ws.Range("E46").Value = '  -4.15%  '
$ws.Range("D47").Value = '27.78'
$ws.Range("E47").Value = '  -3.71%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '2.36'
$ws.Range("E49").Value = '  +1.62%  '
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("D51").Value = '2.58'
$ws.Range("E51").Value = '  +6.54%  '
